$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 274. This shifts the existing rows
# 274..385 down to 275..386 and grows the sheet dimension to A1:R386,
# matching the target diff.
$ws.Rows.Item(274).EntireRow.Insert()

# Populate the newly inserted row 274 with the new data record.
$ws.Range("A274").Value = 3
$ws.Range("B274").Value = "Femacal de La Calera"
$ws.Range("C274").Value = "Coquimbo"
$ws.Range("D274").Value = 44755
$ws.Range("E274").Value = 5
$ws.Range("F274").Value = 100112040
$ws.Range("G274").Value = "Cilantro"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 200
$ws.Range("K274").Value = 3800
$ws.Range("L274").Value = 4000
$ws.Range("M274").Value = 3890
$ws.Range("N274").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O274").Value = "Provincia de Quillota"
$ws.Range("P274").Value = 1297
$ws.Range("Q274").Value = 3
$ws.Range("R274").Value = "Hortaliza"
